$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update first four data columns
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2: B2, D2, E2 cleared; C2 updated
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = -2.9033381510991703
$ws.Range("D2").ClearContents()
$ws.Range("E2").ClearContents()

# Row 3: update all four values
$ws.Range("B3").Value = -3.105531684919832
$ws.Range("C3").Value = 3.9959297561476745
$ws.Range("D3").Value = -0.39966137945635438
$ws.Range("E3").Value = 11.749425093518212

# Update selection to match target
$ws.Range("B1:E3").Select()
